# This script rolls the ENTSOE wind-production dataset forward by one day:
# the oldest day (16.02.2026) is dropped, 17.02.2026 shifts into its place,
# and a freshly fetched day (18.02.2026) is appended at the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Timestamp (A), Notified Production MW (B), Actual Production MW (C)
$data = @(
    @(2, 46070, 532.5549999999999, 268),
    @(3, 46070.01041666666, 543.7329999999999, 0),
    @(4, 46070.02083333334, 550.0650000000001, 248),
    @(5, 46070.03125, 557.7910000000001, 257),
    @(6, 46070.04166666666, 593.66, 288),
    @(7, 46070.05208333334, 619.944, 320),
    @(8, 46070.0625, 640.915, 333),
    @(9, 46070.07291666666, 661.393, 364),
    @(10, 46070.08333333334, 686.08, 396),
    @(11, 46070.09375, 710.361, 391),
    @(12, 46070.10416666666, 732.702, 388),
    @(13, 46070.11458333334, 754.477, 385),
    @(14, 46070.125, 810.939, 403),
    @(15, 46070.13541666666, 831.274, 433),
    @(16, 46070.14583333334, 858.84, 466),
    @(17, 46070.15625, 884.7619999999999, 460),
    @(18, 46070.16666666666, 949.149, 450),
    @(19, 46070.17708333334, 985.63, 470),
    @(20, 46070.1875, 1018.822, 524),
    @(21, 46070.19791666666, 1060.006, 611),
    @(22, 46070.20833333334, 1085.199, 733),
    @(23, 46070.21875, 1133.566, 790),
    @(24, 46070.22916666666, 1177.579, 829),
    @(25, 46070.23958333334, 1254.955, 856),
    @(26, 46070.25, 1275.223, 930),
    @(27, 46070.26041666666, 1316.629, 1044),
    @(28, 46070.27083333334, 1391.646, 1120),
    @(29, 46070.28125, 1432.248, 1215),
    @(30, 46070.29166666666, 1502.424, 1290),
    @(31, 46070.30208333334, 1523.446, 1387),
    @(32, 46070.3125, 1544.784, 1474),
    @(33, 46070.32291666666, 1568.648, 1582),
    @(34, 46070.33333333334, 1579.746, 1617),
    @(35, 46070.34375, 1608.749, 1679),
    @(36, 46070.35416666666, 1631.063, 1712),
    @(37, 46070.36458333334, 1659.334, 1729),
    @(38, 46070.375, 1751.969, 1724),
    @(39, 46070.38541666666, 1776.592, 1683),
    @(40, 46070.39583333334, 1799.963, 1712),
    @(41, 46070.40625, 1824.572, 1727),
    @(42, 46070.41666666666, 1868.847, 1748),
    @(43, 46070.42708333334, 1895.743, 1710),
    @(44, 46070.4375, 1923.719, 1731),
    @(45, 46070.44791666666, 1953.292, 1763),
    @(46, 46070.45833333334, 2007.655, 1820),
    @(47, 46070.46875, 2031.57, 1879),
    @(48, 46070.47916666666, 2053.849, 1940),
    @(49, 46070.48958333334, 2074.652, 1975),
    @(50, 46070.5, 2105.569, 2028),
    @(51, 46070.51041666666, 2128.4, 2068),
    @(52, 46070.52083333334, 2136.944, 2103),
    @(53, 46070.53125, 2161.59, 2151),
    @(54, 46070.54166666666, 2196.402, 2147),
    @(55, 46070.55208333334, 2214.104, 2170),
    @(56, 46070.5625, 2228.06, 2172),
    @(57, 46070.57291666666, 2243.499, 2211),
    @(58, 46070.58333333334, 2285.145, 2320),
    @(59, 46070.59375, 2298.524, 2360),
    @(60, 46070.60416666666, 2311.275, 2451),
    @(61, 46070.61458333334, 2324.355, 2477),
    @(62, 46070.625, 2351.24, 2527),
    @(63, 46070.63541666666, 2359.817, 2551),
    @(64, 46070.64583333334, 2368.477, 2561),
    @(65, 46070.65625, 2310.155, 2507),
    @(66, 46070.66666666666, 2318.99, 2469),
    @(67, 46070.67708333334, 2394.713, 2537),
    @(68, 46070.6875, 2402.908, 2560),
    @(69, 46070.69791666666, 2409.852, 2548),
    @(70, 46070.70833333334, 2429.293, 2518),
    @(71, 46070.71875, 2428.504, 2508),
    @(72, 46070.72916666666, 2428.357, 2504),
    @(73, 46070.73958333334, 2426.213, 2519),
    @(74, 46070.75, 2421.204, 2517),
    @(75, 46070.76041666666, 2419.682, 2529),
    @(76, 46070.77083333334, 2419.217, 2525),
    @(77, 46070.78125, 2418.466, 2501),
    @(78, 46070.79166666666, 2424.162, 2514),
    @(79, 46070.80208333334, 2424.04, 2515),
    @(80, 46070.8125, 2424.333, 2522),
    @(81, 46070.82291666666, 2422.879, 2514),
    @(82, 46070.83333333334, 2417.236, 2521),
    @(83, 46070.84375, 2414.846, 2532),
    @(84, 46070.85416666666, 2412.346, 2522),
    @(85, 46070.86458333334, 2409.994, 2520),
    @(86, 46070.875, 2403.489, 2533),
    @(87, 46070.88541666666, 2334.334, 2517),
    @(88, 46070.89583333334, 2332.247, 2494),
    @(89, 46070.90625, 2331.557, 2470),
    @(90, 46070.91666666666, 2400.096, 2449),
    @(91, 46070.92708333334, 2342.832, 2426),
    @(92, 46070.9375, 2349.07, 2445),
    @(93, 46070.94791666666, 2355.098, 2455),
    @(94, 46070.95833333334, 2379.49, 2432),
    @(95, 46070.96875, 2377.801, 2443),
    @(96, 46070.97916666666, 2376.733, 2421),
    @(97, 46070.98958333334, 2374.961, 2396),
    @(98, 46071, 2463.053, 2354),
    @(99, 46071.01041666666, 2513.594, 2351),
    @(100, 46071.02083333334, 2512.251, 2320),
    @(101, 46071.03125, 2512.578, 2284),
    @(102, 46071.04166666666, 2509.834, 2277),
    @(103, 46071.05208333334, 2509.987, 2253),
    @(104, 46071.0625, 2513.228, 2246),
    @(105, 46071.07291666666, 2513.173, 2250),
    @(106, 46071.08333333334, 2503.831, 2234),
    @(107, 46071.09375, 2498.452, 2236),
    @(108, 46071.10416666666, 2493.971, 2202),
    @(109, 46071.11458333334, 2489.608, 2141),
    @(110, 46071.125, 2486.004, 2095),
    @(111, 46071.13541666666, 2479.585, 2063),
    @(112, 46071.14583333334, 2473.81, 2013),
    @(113, 46071.15625, 2468.688, 1967),
    @(114, 46071.16666666666, 2414.813, 1934),
    @(115, 46071.17708333334, 2396.077, 1892),
    @(116, 46071.1875, 2435.363, 1795),
    @(117, 46071.19791666666, 2425.5, 1745),
    @(118, 46071.20833333334, 2407.195, 1761),
    @(119, 46071.21875, 2396.278, 1734),
    @(120, 46071.22916666666, 2381.012, 1679),
    @(121, 46071.23958333334, 2430.486, 1680),
    @(122, 46071.25, 2346.001, 1706),
    @(123, 46071.26041666666, 2408.635, 1731),
    @(124, 46071.27083333334, 2397.472, 1754),
    @(125, 46071.28125, 2387.696, 1758),
    @(126, 46071.29166666666, 2368.854, 0),
    @(127, 46071.30208333334, 2319.774, 0),
    @(128, 46071.3125, 2282.062, 0),
    @(129, 46071.32291666666, 2275.967, 0),
    @(130, 46071.33333333334, 2266.263, 0),
    @(131, 46071.34375, 2334.075, 0),
    @(132, 46071.35416666666, 2330.932, 0),
    @(133, 46071.36458333334, 2287.097, 0),
    @(134, 46071.375, 2327.294, 0),
    @(135, 46071.38541666666, 2327.533, 0),
    @(136, 46071.39583333334, 2338.478, 0),
    @(137, 46071.40625, 2271.155, 0),
    @(138, 46071.41666666666, 2256.558, 0),
    @(139, 46071.42708333334, 2268.089, 0),
    @(140, 46071.4375, 2266.102, 0),
    @(141, 46071.44791666666, 2263.024, 0),
    @(142, 46071.45833333334, 2263.705, 0),
    @(143, 46071.46875, 2261.263, 0),
    @(144, 46071.47916666666, 2259.501, 0),
    @(145, 46071.48958333334, 2257.677, 0),
    @(146, 46071.5, 2252.241, 0),
    @(147, 46071.51041666666, 2252.011, 0),
    @(148, 46071.52083333334, 2249.037, 0),
    @(149, 46071.53125, 2249.674, 0),
    @(150, 46071.54166666666, 2245.542, 0),
    @(151, 46071.55208333334, 2240.777, 0),
    @(152, 46071.5625, 2236.367, 0),
    @(153, 46071.57291666666, 2231.553, 0),
    @(154, 46071.58333333334, 2234.336, 0),
    @(155, 46071.59375, 2224.506, 0),
    @(156, 46071.60416666666, 2214.764, 0),
    @(157, 46071.61458333334, 2206.038, 0),
    @(158, 46071.625, 2187.782, 0),
    @(159, 46071.63541666666, 2176.793, 0),
    @(160, 46071.64583333334, 2232.889, 0),
    @(161, 46071.65625, 2220.198, 0),
    @(162, 46071.66666666666, 2126.174, 0),
    @(163, 46071.67708333334, 2106.884, 0),
    @(164, 46071.6875, 2152.414, 0),
    @(165, 46071.69791666666, 2131.175, 0),
    @(166, 46071.70833333334, 2114.181, 0),
    @(167, 46071.71875, 2083.805, 0),
    @(168, 46071.72916666666, 2052.356, 0),
    @(169, 46071.73958333334, 2021.27, 0),
    @(170, 46071.75, 1978.173, 0),
    @(171, 46071.76041666666, 1942.199, 0),
    @(172, 46071.77083333334, 1905.566, 0),
    @(173, 46071.78125, 1869.225, 0),
    @(174, 46071.79166666666, 1823.235, 0),
    @(175, 46071.80208333334, 1794.612, 0),
    @(176, 46071.8125, 1765.596, 0),
    @(177, 46071.82291666666, 1737.247, 0),
    @(178, 46071.83333333334, 1682.928, 0),
    @(179, 46071.84375, 1651.05, 0),
    @(180, 46071.85416666666, 1618.691, 0),
    @(181, 46071.86458333334, 1586.118, 0),
    @(182, 46071.875, 1527.538, 0),
    @(183, 46071.88541666666, 1493.436, 0),
    @(184, 46071.89583333334, 1457.91, 0),
    @(185, 46071.90625, 1423.211, 0),
    @(186, 46071.91666666666, 1364.647, 0),
    @(187, 46071.92708333334, 1324.051, 0),
    @(188, 46071.9375, 1283.083, 0),
    @(189, 46071.94791666666, 1243.139, 0),
    @(190, 46071.95833333334, 0, 0),
    @(191, 46071.96875, 0, 0),
    @(192, 46071.97916666666, 0, 0),
    @(193, 46071.98958333334, 0, 0)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]

    # Column E (Lookup) is "<DD>.02.2026<Quarter>"; Quarter (col D) is unchanged (1-96),
    # only the date advances: rows 2-97 -> 17.02.2026, rows 98-193 -> 18.02.2026.
    $quarter = $ws.Cells.Item($r, 4).Value2
    if ($r -le 97) {
        $day = "17"
    } else {
        $day = "18"
    }
    $ws.Cells.Item($r, 5).Value = ($day + ".02.2026" + $quarter)
}

Write-Output "done"